# Data update: refresh "Pagos" (col F) and "Inscrições homologadas" (col H)
# figures in the Inscricoes worksheet (Table1), row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; F=65; H=77},
    @{Row=3; F=21; H=25},
    @{Row=4; F=15; H=29},
    @{Row=5; F=6; H=6},
    @{Row=6; F=6; H=7},
    @{Row=7; F=13; H=14},
    @{Row=8; F=14; H=15},
    @{Row=9; F=24; H=33},
    @{Row=12; F=10; H=10},
    @{Row=15; F=136; H=177},
    @{Row=16; F=7; H=11},
    @{Row=17; F=98; H=130},
    @{Row=18; F=98; H=135},
    @{Row=19; F=65; H=78},
    @{Row=20; F=6; H=9},
    @{Row=21; F=3; H=3},
    @{Row=23; F=7; H=9},
    @{Row=24; F=26; H=30},
    @{Row=25; F=22; H=30},
    @{Row=26; F=31; H=41},
    @{Row=27; F=18; H=22},
    @{Row=28; F=25; H=27},
    @{Row=29; F=16; H=19},
    @{Row=30; F=3; H=3},
    @{Row=32; F=10; H=19},
    @{Row=33; F=24; H=36},
    @{Row=34; F=21; H=24},
    @{Row=36; F=87; H=119},
    @{Row=37; F=58; H=70},
    @{Row=38; F=42; H=62},
    @{Row=39; F=21; H=29},
    @{Row=40; F=24; H=26},
    @{Row=41; F=34; H=45},
    @{Row=42; F=36; H=45},
    @{Row=43; F=29; H=32},
    @{Row=44; F=25; H=35},
    @{Row=45; F=20; H=27},
    @{Row=46; F=18; H=27},
    @{Row=47; F=60; H=70},
    @{Row=48; F=42; H=48},
    @{Row=49; F=55; H=72},
    @{Row=50; F=17; H=26},
    @{Row=51; F=13; H=13},
    @{Row=52; F=10; H=10},
    @{Row=55; F=8; H=11},
    @{Row=56; F=8; H=10},
    @{Row=57; F=17; H=21},
    @{Row=60; F=15; H=20},
    @{Row=61; F=25; H=35},
    @{Row=62; F=28; H=42},
    @{Row=63; F=26; H=34},
    @{Row=64; F=30; H=35},
    @{Row=65; F=26; H=39},
    @{Row=67; F=35; H=43},
    @{Row=68; F=13; H=17},
    @{Row=69; F=15; H=18},
    @{Row=70; F=38; H=51},
    @{Row=71; F=36; H=46},
    @{Row=72; F=41; H=52},
    @{Row=73; F=25; H=37},
    @{Row=74; F=11; H=15},
    @{Row=75; F=14; H=19},
    @{Row=76; F=29; H=46},
    @{Row=77; F=41; H=58},
    @{Row=79; F=39; H=50},
    @{Row=81; F=18; H=23},
    @{Row=82; F=12; H=18},
    @{Row=83; F=7; H=14},
    @{Row=86; F=4; H=6},
    @{Row=87; F=14; H=21},
    @{Row=88; F=27; H=35},
    @{Row=89; F=29; H=35}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    $ws.Cells.Item($u.Row, 8).Value = $u.H
}
